# "Add event when regedit open"
# Adds a new localization entry (key "reg_1"/"reg_2") describing an event
# that fires when the player opens the Windows registry editor, to the
# keys / ru-RU / en-US sheets of the locales workbook.

$wb = $excel.ActiveWorkbook

$keysSheet = $wb.Worksheets.Item("keys")
$ruSheet   = $wb.Worksheets.Item("ru-RU")
$enSheet   = $wb.Worksheets.Item("en-US")

# --- keys sheet: new key names + locale count (row 32/33) -------------
$keysSheet.Range("A32").Value = "reg_1"
$keysSheet.Range("B32").Value = 4
$keysSheet.Range("A33").Value = "reg_2"
$keysSheet.Range("B33").Value = 4

# --- ru-RU translations -------------------------------------------------
$ruSheet.Range("A32").Value = "reg_1"
$ruSheet.Range("B32").Value = "Если честно, то я уже сделала пару записей в реестре."
$ruSheet.Range("A33").Value = "reg_2"
$ruSheet.Range("B33").Value = "Поищи их, если хочешь."

# --- en-US translations --------------------------------------------------
$enSheet.Range("A32").Value = "reg_1"
$enSheet.Range("B32").Value = "To be honest, I have already made a couple entries in the registry."
$enSheet.Range("A33").Value = "reg_2"
$enSheet.Range("B33").Value = "If you want, try to find and look for them."

# Widen the translation column on en-US now that it holds a longer string.
$enSheet.Columns.Item(2).ColumnWidth = 63.166666666666664

# --- restore/update each sheet's on-screen selection ---------------------
$keysSheet.Range("A32:B33").Select()
$ruSheet.Range("B33").Select()
$enSheet.Range("B33").Select()

# en-US becomes the active (visible) tab when the workbook is saved.
$enSheet.Activate()
